# Update the three question/answer cells whose linked URLs were rotated
# to new endpoints, while leaving their existing hyperlink targets intact
# (only the displayed text changes, matching the upstream edit).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C14").Value = "https://raw.githubusercontent.com/Telvinvarghese/test/main/email.json"
$ws.Range("C26").Value = "https://github.com/Telvinvarghese/test"
$ws.Range("C29").Value = "https://3db5-223-178-84-140.ngrok-free.app/"

# Restore the selection/scroll position recorded in the saved workbook.
$ws.Range("A21").Select()
$ws.Range("B16").Select()
